$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 2).Value = "BBB"
$ws.Cells.Item(9, 2).Value = "A"
$ws.Cells.Item(10, 2).Value = "A"
$ws.Cells.Item(11, 2).Value = "A"
$ws.Cells.Item(12, 2).Value = "AAA"
$ws.Cells.Item(18, 2).Value = "BBB"
$ws.Cells.Item(19, 2).Value = "BBB"
$ws.Cells.Item(22, 2).Value = "B"
$ws.Cells.Item(28, 2).Value = "A"
$ws.Cells.Item(29, 2).Value = "BB"
$ws.Cells.Item(30, 2).Value = "A"
$ws.Cells.Item(33, 2).Value = "BBB"
$ws.Cells.Item(36, 2).Value = "BBB"
$ws.Cells.Item(42, 2).Value = "BBB"
$ws.Cells.Item(43, 2).Value = "BBB"
$ws.Cells.Item(45, 2).Value = "BBB"
$ws.Cells.Item(47, 2).Value = "BB"
$ws.Cells.Item(51, 2).Value = "BB"
$ws.Cells.Item(55, 2).Value = "BB"
$ws.Cells.Item(58, 2).Value = "BB"
$ws.Cells.Item(59, 2).Value = "BB"
$ws.Cells.Item(60, 2).Value = "BB"
$ws.Cells.Item(61, 2).Value = "A"
$ws.Cells.Item(63, 2).Value = "A"
$ws.Cells.Item(65, 2).Value = "BBB"
$ws.Cells.Item(66, 2).Value = "BBB"
$ws.Cells.Item(68, 2).Value = "BB"
$ws.Cells.Item(69, 2).Value = "BB"
$ws.Cells.Item(70, 2).Value = "BB"
$ws.Cells.Item(71, 2).Value = "BBB"
$ws.Cells.Item(73, 2).Value = "A"
$ws.Cells.Item(80, 2).Value = "B"
$ws.Cells.Item(83, 2).Value = "AA"
$ws.Cells.Item(85, 2).Value = "BB"
$ws.Cells.Item(88, 2).Value = "BB"
$ws.Cells.Item(90, 2).Value = "BBB"
$ws.Cells.Item(91, 2).Value = "A"
$ws.Cells.Item(92, 2).Value = "A"
$ws.Cells.Item(95, 2).Value = "BBB"
$ws.Cells.Item(99, 2).Value = "A"
$ws.Cells.Item(103, 2).Value = "BBB"
$ws.Cells.Item(110, 2).Value = "BBB"
$ws.Cells.Item(112, 2).Value = "A"
$ws.Cells.Item(113, 2).Value = "A"
$ws.Cells.Item(120, 2).Value = "A"
$ws.Cells.Item(122, 2).Value = "BB"
$ws.Cells.Item(123, 2).Value = "BBB"
$ws.Cells.Item(124, 2).Value = "A"
$ws.Cells.Item(125, 2).Value = "A"
$ws.Cells.Item(127, 2).Value = "BBB"
$ws.Cells.Item(130, 2).Value = "A"
$ws.Cells.Item(131, 2).Value = "BB"
$ws.Cells.Item(134, 2).Value = "A"
$ws.Cells.Item(135, 2).Value = "BBB"
$ws.Cells.Item(138, 2).Value = "BBB"
$ws.Cells.Item(140, 2).Value = "BBB"
$ws.Cells.Item(142, 2).Value = "A"
$ws.Cells.Item(144, 2).Value = "A"
$ws.Cells.Item(145, 2).Value = "A"
$ws.Cells.Item(148, 2).Value = "BBB"
$ws.Cells.Item(150, 2).Value = "BB"
$ws.Cells.Item(152, 2).Value = "B"
$ws.Cells.Item(154, 2).Value = "BBB"
$ws.Cells.Item(155, 2).Value = "BB"
$ws.Cells.Item(156, 2).Value = "B"
$ws.Cells.Item(158, 2).Value = "BBB"
$ws.Cells.Item(159, 2).Value = "AA"
$ws.Cells.Item(160, 2).Value = "A"
$ws.Cells.Item(162, 2).Value = "A"
$ws.Cells.Item(166, 2).Value = "BBB"
$ws.Cells.Item(169, 2).Value = "BBB"
$ws.Cells.Item(171, 2).Value = "BB"
$ws.Cells.Item(173, 2).Value = "BBB"
$ws.Cells.Item(178, 2).Value = "A"
$ws.Cells.Item(182, 2).Value = "BBB"
$ws.Cells.Item(185, 2).Value = "BBB"
$ws.Cells.Item(186, 2).Value = "AAA"
$ws.Cells.Item(187, 2).Value = "AAA"
$ws.Cells.Item(196, 2).Value = "BB"
$ws.Cells.Item(204, 2).Value = "A"
$ws.Cells.Item(205, 2).Value = "AAA"
$ws.Cells.Item(206, 2).Value = "BB"
$ws.Cells.Item(208, 2).Value = "BBB"
$ws.Cells.Item(213, 2).Value = "B"
$ws.Cells.Item(214, 2).Value = "BBB"
$ws.Cells.Item(215, 2).Value = "BB"
$ws.Cells.Item(220, 2).Value = "BBB"
$ws.Cells.Item(227, 2).Value = "BBB"
$ws.Cells.Item(228, 2).Value = "B"
$ws.Cells.Item(230, 2).Value = "BBB"
$ws.Cells.Item(231, 2).Value = "BB"
$ws.Cells.Item(235, 2).Value = "BBB"
$ws.Cells.Item(238, 2).Value = "BB"
$ws.Cells.Item(239, 2).Value = "AAA"
$ws.Cells.Item(241, 2).Value = "AAA"
$ws.Cells.Item(244, 2).Value = "BBB"
$ws.Cells.Item(245, 2).Value = "BBB"
$ws.Cells.Item(247, 2).Value = "BBB"
$ws.Cells.Item(252, 2).Value = "BB"
$ws.Cells.Item(256, 2).Value = "BBB"
$ws.Cells.Item(257, 2).Value = "BB"
$ws.Cells.Item(258, 2).Value = "BBB"
$ws.Cells.Item(262, 2).Value = "A"
$ws.Cells.Item(269, 2).Value = "BBB"
$ws.Cells.Item(271, 2).Value = "BB"
$ws.Cells.Item(274, 2).Value = "BB"
$ws.Cells.Item(282, 2).Value = "A"
$ws.Cells.Item(283, 2).Value = "BBB"
$ws.Cells.Item(288, 2).Value = "A"
$ws.Cells.Item(289, 2).Value = "B"
$ws.Cells.Item(292, 2).Value = "A"
$ws.Cells.Item(301, 2).Value = "A"
$ws.Cells.Item(306, 2).Value = "BBB"
$ws.Cells.Item(309, 2).Value = "BBB"
$ws.Cells.Item(313, 2).Value = "BBB"
$ws.Cells.Item(316, 2).Value = "BBB"
$ws.Cells.Item(319, 2).Value = "BBB"
$ws.Cells.Item(322, 2).Value = "BB"
$ws.Cells.Item(330, 2).Value = "A"
$ws.Cells.Item(331, 2).Value = "BBB"
$ws.Cells.Item(332, 2).Value = "B"
$ws.Cells.Item(334, 2).Value = "BB"
$ws.Cells.Item(335, 2).Value = "BB"
$ws.Cells.Item(336, 2).Value = "BB"
$ws.Cells.Item(337, 2).Value = "BB"
$ws.Cells.Item(344, 2).Value = "BBB"
$ws.Cells.Item(346, 2).Value = "BBB"
$ws.Cells.Item(349, 2).Value = "BB"
$ws.Cells.Item(350, 2).Value = "B"
$ws.Cells.Item(356, 2).Value = "A"
$ws.Cells.Item(357, 2).Value = "BBB"
$ws.Cells.Item(361, 2).Value = "A"
$ws.Cells.Item(363, 2).Value = "BB"
$ws.Cells.Item(365, 2).Value = "B"
$ws.Cells.Item(366, 2).Value = "A"
$ws.Cells.Item(370, 2).Value = "BB"
$ws.Cells.Item(371, 2).Value = "BB"
$ws.Cells.Item(373, 2).Value = "BB"
$ws.Cells.Item(376, 2).Value = "A"
$ws.Cells.Item(380, 2).Value = "B"
$ws.Cells.Item(382, 2).Value = "B"
$ws.Cells.Item(383, 2).Value = "A"
$ws.Cells.Item(384, 2).Value = "BB"
$ws.Cells.Item(385, 2).Value = "BB"
$ws.Cells.Item(394, 2).Value = "BBB"
$ws.Cells.Item(396, 2).Value = "B"
$ws.Cells.Item(397, 2).Value = "BBB"
$ws.Cells.Item(399, 2).Value = "BBB"
$ws.Cells.Item(400, 2).Value = "BBB"
